$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated K column (column G) values for rows 2-13
$kValues = @{
    2  = 1
    3  = 0
    4  = 5
    5  = 2
    6  = 2
    7  = 0
    8  = 3
    9  = 1
    10 = 2
    11 = 1
    12 = 2
    13 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
